$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the text of the rows we're about to move before we overwrite anything.
$row4B = $ws.Range("B4").Value2
$row4C = $ws.Range("C4").Value2
$row4D = $ws.Range("D4").Value2

$row9B = $ws.Range("B9").Value2
$row9C = $ws.Range("C9").Value2
$row9D = $ws.Range("D9").Value2

$row10B = $ws.Range("B10").Value2
$row10D = $ws.Range("D10").Value2

# Move the two still-outstanding features down below the separator row,
# into rows 14 and 15 (matching their original, unstyled formatting).
$ws.Range("B14").Value2 = $row4B
$ws.Range("C14").Value2 = $row4C
$ws.Range("D14").Value2 = $row4D

$ws.Range("B15").Value2 = $row9B
$ws.Range("C15").Value2 = $row9C
$ws.Range("D15").Value2 = $row9D

# Clear out the old row 9 and row 10 content/styling - it has moved.
$ws.Range("B9:D9").Clear()
$ws.Range("B10:D10").Clear()

# The "send ball to set player" feature (old row 10) is now finished - move
# it up to the top of the completed list (row 4) and give it a developer
# plus the same strikethrough styling as the rest of the completed rows.
$ws.Range("B4:D4").Value2 = $null
$ws.Range("B4").Value2 = $row10B
$ws.Range("C4").Value2 = "Garrett"
$ws.Range("D4").Value2 = $row10D
$ws.Range("B4:D4").Font.Strikethrough = $true

$ws.Range("D12").Select() | Out-Null
